$d = $word.ActiveDocument

# Locate the paragraph "LOB1036: Geometria Analítica (Requisito fraco)".
# It is kept as-is. Immediately after it the document currently has four
# paragraphs that need to be removed:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) an empty paragraph
#   4) an empty paragraph with PageBreakBefore
# After the deletion, the paragraph right after LOB1036 becomes the empty
# paragraph (originally the 6th paragraph below LOB1036) followed by the
# empty page-break paragraph that used to close the section.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*LOB1036*Geometria*Anal*tica*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Deleting the paragraph right after the target, four times in a row,
    # removes the four unwanted paragraphs while leaving LOB1036's
    # paragraph (and everything beyond the deleted block) untouched.
    for ($k = 0; $k -lt 4; $k++) {
        $d.Paragraphs($targetIndex + 1).Range.Delete()
    }
}
